$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34-75 down to 35-76.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly record.
$ws.Range("A34").Value = 11
$ws.Range("B34").Value = 'Vega Monumental Concepción'
$ws.Range("C34").Value = 'Bíobío'
$ws.Range("D34").Value = (Get-Date -Year 2021 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E34").Value = 8
$ws.Range("F34").Value = 100112043
$ws.Range("G34").Value = 'Pepino ensalada'
$ws.Range("H34").Value = 'Sin especificar'
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 350
$ws.Range("K34").Value = 10000
$ws.Range("L34").Value = 11000
$ws.Range("M34").Value = 10429
$ws.Range("N34").Value = '$/caja 60 unidades'
$ws.Range("O34").Value = 'Región de Arica y Parinacota'
$ws.Range("P34").Value = 174
$ws.Range("Q34").Value = 60
$ws.Range("R34").Value = 'Hortaliza'
